$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right 6 -> 9, Wrong 3 -> 2
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Row 12 ("Total"): Right 84 -> 126, Wrong -6 -> -4, Max fraction 78/168 -> 122/252
$ws.Range("B12").Value = 126
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "122/252"
